# Updates cryptos list prices (column D) and 1h volume % (column E)
# on Sheet1, matching the "Updated cryptos list ... with GitHub Actions" commit.
#
# Column D values are plain numeric-looking text (e.g. "58.53", "1.00") stored
# as strings in the source sheet. Writing them via a plain .Value assignment
# would let Excel auto-coerce them into real numbers, which is not what the
# original workbook does. To keep them as text (and to avoid leaving a lasting
# NumberFormat/style change on the cell), each D cell is briefly switched to a
# text format, given its new literal value, and then returned to the "Normal"
# cell style so the saved style index for the cell is unchanged.
#
# Column E values already contain a "%" sign wrapped in two leading/trailing
# spaces, so they are never ambiguous with numbers and can be set directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.797.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.19%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.077.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.56%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.94%  "

# Row 6
$ws.Range("E6").Value = "  -0.17%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.53"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.48%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("E9").Value = "  +0.21%  "

# Row 10
$ws.Range("E10").Value = "  -0.86%  "

# Row 11
$ws.Range("E11").Value = "  +3.40%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.383.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.49%  "

# Row 13
$ws.Range("E13").Value = "  +0.20%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.49%  "

# Row 15
$ws.Range("E15").Value = "  +1.06%  "

# Row 16
$ws.Range("E16").Value = "  +0.70%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.081.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.52%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.695.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.27%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.94%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.31%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0842"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.35%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.19%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.49%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.09%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.51%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.138"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.98%  "

# Row 29
$ws.Range("E29").Value = "  -1.95%  "

# Row 30
$ws.Range("E30").Value = "  -0.95%  "

# Row 31
$ws.Range("E31").Value = "  +1.10%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.40%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0633"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.14%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.01%  "

# Row 35
$ws.Range("E35").Value = "  -2.56%  "

# Row 36
$ws.Range("E36").Value = "  -0.01%  "

# Row 37
$ws.Range("E37").Value = "  -2.85%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.02%  "

# Row 39
$ws.Range("E39").Value = "  -1.19%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0232"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.44%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.65%  "

# Row 42
$ws.Range("E42").Value = "  -1.14%  "

# Row 43
$ws.Range("E43").Value = "  -0.64%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.87%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.451.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.95%  "

# Row 46
$ws.Range("E46").Value = "  -1.87%  "

# Row 47
$ws.Range("E47").Value = "  -0.96%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.09%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.08%  "

# Row 50
$ws.Range("E50").Value = "  -1.49%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.269.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.47%  "
